$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '60.682.24'
$ws.Range('E2').Value = '  -2.65%  '
$ws.Range('D3').Value = '2.903.16'
$ws.Range('E3').Value = '  -3.80%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'586.08"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').Value = "'147.13"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.31%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -2.60%  '
$ws.Range('D9').Value = '2.903.25'
$ws.Range('E9').Value = '  -3.77%  '
$ws.Range('E10').Value = '  +4.61%  '
$ws.Range('E11').Value = '  -4.14%  '
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('E13').Value = '  -3.52%  '
$ws.Range('D14').Value = "'33.98"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '3.385.76'
$ws.Range('E16').Value = '  -3.76%  '
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').Value = '60.586.90'
$ws.Range('E18').Value = '  -2.74%  '
$ws.Range('D19').Value = '2.904.33'
$ws.Range('E19').Value = '  -3.94%  '
$ws.Range('D20').Value = "'427.80"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.60%  '
$ws.Range('D21').Value = "'13.62"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.20%  '
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').Value = "'7.09"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.39%  '
$ws.Range('D25').Value = "'11.07"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.01%  '
$ws.Range('D26').Value = "'2.20"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.01%  '
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').Value = "'7.19"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('E32').Value = '  -3.26%  '
$ws.Range('E33').Value = '  -3.92%  '
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').Value = '0.0₃0836'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = "'49.29"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = "'2.03"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('D40').Value = "'2.95"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').Value = "'0.123"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = "'8.74"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').Value = "'41.90"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').Value = "'369.92"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.93%  '
$ws.Range('D47').Value = "'133.80"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = '2.654.98'
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('D50').Value = "'24.94"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.12%  '
$ws.Range('E51').Value = '  -1.25%  '
